$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 340 (shifts existing rows 340-363 down to 341-364,
# carrying their formatting/styles along, same as Excel's native row insert).
$ws.Rows.Item(340).Insert()

# Populate the newly inserted row 340 with a new weekly price record.
$ws.Range("A340").Value = 7
$ws.Range("B340").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C340").Value = "Ñuble"
$ws.Range("D340").Value = 45021
$ws.Range("E340").Value = 16
$ws.Range("F340").Value = 100112006
$ws.Range("G340").Value = "Repollo"
$ws.Range("H340").Value = "Crespo record"
$ws.Range("I340").Value = "Primera"
$ws.Range("J340").Value = 250
$ws.Range("K340").Value = 1300
$ws.Range("L340").Value = 1300
$ws.Range("M340").Value = 1300
$ws.Range("N340").Value = "$/unidad"
$ws.Range("O340").Value = "Provincia de Diguillín"
$ws.Range("P340").Value = 1300
$ws.Range("Q340").Value = 1
$ws.Range("R340").Value = "Hortaliza"
